$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5203476666666668
$ws.Range("N2").Value = 1.561043
$ws.Range("O2").Value = 0.004105934376266647
$ws.Range("P2").Value = 0.004105934376266647
$ws.Range("Q2").Value = 0.08516443535722223
$ws.Range("R2").Value = 0.7664799182150002
$ws.Range("S2").Value = 0.00003547319153750749
$ws.Range("T2").Value = 0.0000354731915375075
$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("O3").Value = 0.8361295370252257
$ws.Range("P3").Value = 0.8361295370252259
$ws.Range("Q3").Value = 17.34282464859945
$ws.Range("R3").Value = 156.085421837395
$ws.Range("S3").Value = 0.007223735330137461
$ws.Range("T3").Value = 0.007223735330137464
$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 20.24706
$ws.Range("N4").Value = 60.74118
$ws.Range("O4").Value = 0.1597645285985076
$ws.Range("P4").Value = 0.1597645285985076
$ws.Range("Q4").Value = 3.3138025651
$ws.Range("R4").Value = 29.8242230859
$ws.Range("S4").Value = 0.001380284535630485
$ws.Range("T4").Value = 0.001380284535630486
$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5203476666666668
$ws.Range("N5").Value = 1.561043
$ws.Range("O5").Value = 0.004105934376266647
$ws.Range("P5").Value = 0.004105934376266647
$ws.Range("Q5").Value = 7.973196755839336
$ws.Range("R5").Value = 71.75877080255401
$ws.Range("S5").Value = 0.003321042809710083
$ws.Range("T5").Value = 0.003321042809710084
$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("O6").Value = 0.8361295370252257
$ws.Range("P6").Value = 0.8361295370252259
$ws.Range("S6").Value = 0.6762947803000919
$ws.Range("T6").Value = 0.6762947803000922
$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 20.24706
$ws.Range("N7").Value = 60.74118
$ws.Range("O7").Value = 0.1597645285985076
$ws.Range("P7").Value = 0.1597645285985076
$ws.Range("Q7").Value = 310.2421773915601
$ws.Range("R7").Value = 2792.17959652404
$ws.Range("S7").Value = 0.129223896518101
$ws.Range("T7").Value = 0.129223896518101
$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5203476666666668
$ws.Range("N8").Value = 1.561043
$ws.Range("O8").Value = 0.004105934376266647
$ws.Range("P8").Value = 0.004105934376266647
$ws.Range("Q8").Value = 1.799212024306889
$ws.Range("R8").Value = 16.192908218762
$ws.Range("S8").Value = 0.0007494183750190555
$ws.Range("T8").Value = 0.0007494183750190556
$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("O9").Value = 0.8361295370252257
$ws.Range("P9").Value = 0.8361295370252259
$ws.Range("Q9").Value = 366.3902486093318
$ws.Range("R9").Value = 3297.512237483986
$ws.Range("S9").Value = 0.1526110213949962
$ws.Range("T9").Value = 0.1526110213949963
$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 20.24706
$ws.Range("N10").Value = 60.74118
$ws.Range("O10").Value = 0.1597645285985076
$ws.Range("P10").Value = 0.1597645285985076
$ws.Range("Q10").Value = 70.00848882868002
$ws.Range("R10").Value = 630.07639945812
$ws.Range("S10").Value = 0.02916034754477612
$ws.Range("T10").Value = 0.02916034754477612

Write-Host "Updated cells"
